{"js": "const replacements = [\n  [\"812\u00f74=\", \"535\u00f73=\"],\n  [\"681\u00f79=\", \"364\u00f72=\"],\n  [\"123\u00f75=\", \"564\u00f79=\"],\n  [\"896\u00f76=\", \"731\u00f73=\"],\n  [\"156\u00f76=\", \"289\u00f72=\"],\n  [\"716\u00f75=\", \"936\u00f75=\"],\n  [\"217\u00f79=\", \"973\u00f73=\"],\n  [\"285\u00f73=\", \"485\u00f72=\"],\n  [\"841\u00f72=\", \"840\u00f79=\"],\n  [\"437\u00f74=\", \"427\u00f78=\"],\n  [\"418\u00f75=\", \"282\u00f76=\"],\n  [\"533\u00f79=\", \"135\u00f73=\"],\n  [\"832\u00f79=\", \"190\u00f78=\"],\n  [\"832\u00f74=\", \"852\u00f78=\"],\n  [\"341\u00f78=\", \"440\u00f73=\"],\n  [\"198\u00f73=\", \"813\u00f79=\"],\n  [\"410\u00f72=\", \"425\u00f74=\"],\n  [\"288\u00f79=\", \"803\u00f79=\"],\n  [\"281\u00f76=\", \"842\u00f72=\"],\n  [\"484\u00f78=\", \"586\u00f77=\"],\n  [\"847\u00f74=\", \"647\u00f75=\"],\n  [\"941\u00f76=\", \"482\u00f78=\"],\n  [\"725\u00f75=\", \"747\u00f75=\"],\n  [\"740\u00f79=\", \"648\u00f77=\"],\n  [\"736\u00f76=\", \"310\u00f75=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old = \"812\u00f74=\"; New = \"535\u00f73=\"},\n    @{Old = \"681\u00f79=\"; New = \"364\u00f72=\"},\n    @{Old = \"123\u00f75=\"; New = \"564\u00f79=\"},\n    @{Old = \"896\u00f76=\"; New = \"731\u00f73=\"},\n    @{Old = \"156\u00f76=\"; New = \"289\u00f72=\"},\n    @{Old = \"716\u00f75=\"; New = \"936\u00f75=\"},\n    @{Old = \"217\u00f79=\"; New = \"973\u00f73=\"},\n    @{Old = \"285\u00f73=\"; New = \"485\u00f72=\"},\n    @{Old = \"841\u00f72=\"; New = \"840\u00f79=\"},\n    @{Old = \"437\u00f74=\"; New = \"427\u00f78=\"},\n    @{Old = \"418\u00f75=\"; New = \"282\u00f76=\"},\n    @{Old = \"533\u00f79=\"; New = \"135\u00f73=\"},\n    @{Old = \"832\u00f79=\"; New = \"190\u00f78=\"},\n    @{Old = \"832\u00f74=\"; New = \"852\u00f78=\"},\n    @{Old = \"341\u00f78=\"; New = \"440\u00f73=\"},\n    @{Old = \"198\u00f73=\"; New = \"813\u00f79=\"},\n    @{Old = \"410\u00f72=\"; New = \"425\u00f74=\"},\n    @{Old = \"288\u00f79=\"; New = \"803\u00f79=\"},\n    @{Old = \"281\u00f76=\"; New = \"842\u00f72=\"},\n    @{Old = \"484\u00f78=\"; New = \"586\u00f77=\"},\n    @{Old = \"847\u00f74=\"; New = \"647\u00f75=\"},\n    @{Old = \"941\u00f76=\"; New = \"482\u00f78=\"},\n    @{Old = \"725\u00f75=\"; New = \"747\u00f75=\"},\n    @{Old = \"740\u00f79=\"; New = \"648\u00f77=\"},\n    @{Old = \"736\u00f76=\"; New = \"310\u00f75=\"}\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $r.New\n    $find.Execute(\n        $r.Old,      # FindText\n        $false,      # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $r.New,      # ReplaceWith\n        2            # Replace (wdReplaceAll)\n    )\n}\n"}
